$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
for ($i=1; $i -le $sm.CustomLayouts.Count; $i++) {
    $lay = $sm.CustomLayouts.Item($i)
    for ($j=1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "10/25/2019"
        }
    }
}
Write-Output "done"
